# ironing, cleaning up code
#
# Updates the BYTE/ORIGINAL/REBUILT/NOTES comparison table:
#  - fill in a missing REBUILT value that mirrors ORIGINAL (byte 12 row)
#  - rename the "ASCII new page" note to "size of data subblock"
#  - add the missing "img!" note for the byte-23 row
#  - move the selection/viewport to where we left off editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 (BYTE 12): REBUILT column was left blank -- mirror ORIGINAL (255)
$ws.Range("C20").Value = 255

# Row 23 (BYTE 15): NOTES text update
$ws.Range("D23").Value = "size of data subblock"

# Row 37 (BYTE 23): give it its missing NOTES label
$ws.Range("D37").Value = "img!"

# Move viewport/selection to reflect where we left off editing
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("E21").Select()
